$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit removed the post that used to occupy row 435 ("「見よ」...").
# All the rows below it (436-587) shift up by one to fill the gap, which is
# exactly what EntireRow delete-with-shift-up does.
$ws.Rows.Item(435).Delete()
